# Refresh the "cryptos" price/volume snapshot (Price = column D, Volume(1h) = column E).
# Price cells that look numeric (e.g. "316.87", "0.09140") are forced to Text format
# first so Excel doesn't silently coerce them to numbers and drop significant/trailing
# digits; the "xx.xxx.xx" style prices (e.g. "24.605.77") are already non-numeric text
# and don't need that treatment.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.605.77"
$ws.Range("E2").Value = "  +3.38%  "
$ws.Range("D3").Value = "1.695.71"
$ws.Range("E3").Value = "  +2.05%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.87"
$ws.Range("E5").Value = "  +2.32%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("E7").Value = "  +1.75%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4018"
$ws.Range("E8").Value = "  +1.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.528"
$ws.Range("E9").Value = "  +4.95%  "
$ws.Range("E10").Value = "  -0.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.67"
$ws.Range("E11").Value = "  -0.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08768"
$ws.Range("E12").Value = "  +1.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.219"
$ws.Range("E13").Value = "  +6.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.29"
$ws.Range("E14").Value = "  +3.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.206"
$ws.Range("E15").Value = "  +12.43%  "
$ws.Range("E16").Value = "  +0.79%  "
$ws.Range("D17").Value = "1.698.23"
$ws.Range("E17").Value = "  +2.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "99.80"
$ws.Range("E18").Value = "  +0.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07073"
$ws.Range("E19").Value = "  +2.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.69"
$ws.Range("E20").Value = "  +3.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.965"
$ws.Range("E21").Value = "  +4.95%  "
$ws.Range("E22").Value = "  -0.38%  "
$ws.Range("E23").Value = "  +3.19%  "
$ws.Range("D24").Value = "24.610.40"
$ws.Range("E24").Value = "  +3.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.127"
$ws.Range("E25").Value = "  +9.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.334"
$ws.Range("E26").Value = "  +0.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.87"
$ws.Range("E27").Value = "  +5.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.76"
$ws.Range("E28").Value = "  +1.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "136.16"
$ws.Range("E29").Value = "  +4.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.187"
$ws.Range("E30").Value = "  +1.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.486"
$ws.Range("E31").Value = "  +10.08%  "
$ws.Range("D32").Value = "1.881.43"
$ws.Range("E32").Value = "  +1.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.090"
$ws.Range("E33").Value = "  -1.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08612"
$ws.Range("E34").Value = "  +1.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.142"
$ws.Range("E35").Value = "  +7.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.63"
$ws.Range("E36").Value = "  +11.29%  "
$ws.Range("E37").Value = "  +3.59%  "
$ws.Range("E38").Value = "  +0.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.48"
$ws.Range("E39").Value = "  +0.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09140"
$ws.Range("E40").Value = "  +4.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02743"
$ws.Range("E41").Value = "  +8.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.483"
$ws.Range("E42").Value = "  +2.13%  "
$ws.Range("E43").Value = "  +1.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7173"
$ws.Range("E44").Value = "  +1.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.60"
$ws.Range("E45").Value = "  +4.66%  "
$ws.Range("E46").Value = "  +5.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.222"
$ws.Range("E47").Value = "  +2.67%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.01"
$ws.Range("E49").Value = "  +1.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.319"
$ws.Range("E50").Value = "  +8.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07985"
$ws.Range("E51").Value = "  +2.52%  "
